# Update the "想去人数" (column F) figures on the "展览" and "全部类型" sheets
# to match the refreshed output data (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Column F = 6 on every affected sheet.
$col = 6

# Row => new value, for sheet "展览"
$exhibitionUpdates = @{
    2  = 859
    5  = 1176
    6  = 56
    7  = 4244
    8  = 2566
    10 = 2441
    14 = 1639
    15 = 652
    16 = 12
    17 = 105
    18 = 310
    20 = 267
    22 = 5
    23 = 457
    25 = 86
    26 = 504
    27 = 680
    28 = 91
    30 = 381
    33 = 923
    34 = 65
    36 = 1001
    37 = 1979
    38 = 239
    39 = 4
    40 = 523
    41 = 79
    42 = 12
    43 = 611
    44 = 1274
    45 = 63
    47 = 415
    48 = 58
}

# Row => new value, for sheet "全部类型"
$allTypesUpdates = @{
    2  = 859
    3  = 1176
    5  = 56
    6  = 4244
    7  = 2566
    8  = 2441
    10 = 1639
    12 = 652
    13 = 12
    14 = 105
    15 = 310
    17 = 267
    19 = 457
    21 = 86
    22 = 504
    23 = 680
    24 = 91
    29 = 381
    31 = 923
    32 = 65
    35 = 1001
    36 = 1979
    37 = 239
    42 = 79
    43 = 12
    44 = 611
    45 = 1274
    46 = 63
    47 = 415
    48 = 58
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $ws1.Cells.Item($row, $col).Value = $exhibitionUpdates[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $ws4.Cells.Item($row, $col).Value = $allTypesUpdates[$row]
}
